# Update cryptos list values (prices and 1h volume deltas) to match latest scrape.
# Also swaps rows 43/44 (WEMIXToken/Aave) order with updated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.910.43"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "'1.666.10"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'215.48"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'0.533"
$ws.Range("E6").Value = "  +4.95%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'20.23"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("D12").Value = "'1.901.48"
$ws.Range("D13").Value = "'1.652.81"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "'0.524"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'66.04"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "'26.926.73"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "'234.80"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "'7.99"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "'0.0₃0732"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'4.34"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").Value = "'2.20"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").Value = "'9.09"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "'145.97"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "'7.11"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").Value = "'15.87"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").Value = "'1.455.42"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'0.580"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'0.903"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'5.71"
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'0.973"
$ws.Range("E43").Value = "  +6.06%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'65.78"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'1.809.17"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "'90.55"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'0.0₆0104"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  +4.34%  "
$ws.Range("E51").Value = "  +0.03%  "
